# Weekly refresh of the Espinaca (Vega Central Mapocho de Santiago) price series.
# Two new observations (most recent week) are inserted at the top of the
# chronologically-ordered block (rows 209-210); every existing row from 209-249
# shifts down by two rows (-> 211-251); the two rows pushed past the former last
# row (249) are appended as the new rows 250-251.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 209
$ws.Cells.Item(209, 1).Value = 9  # Mercado ID
$ws.Cells.Item(209, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(209, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(209, 4).Value = 44476  # Fecha
$ws.Cells.Item(209, 5).Value = 13  # Codreg
$ws.Cells.Item(209, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(209, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(209, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(209, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(209, 10).Value = 250  # Volumen
$ws.Cells.Item(209, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(209, 12).Value = 7000  # Precio máximo
$ws.Cells.Item(209, 13).Value = 6500  # Precio promedio ponderado
$ws.Cells.Item(209, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(209, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(209, 16).Value = 650  # Precio $/Kg
$ws.Cells.Item(209, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(209, 18).Value = 'Hortaliza'  # Clasificación

# Row 210
$ws.Cells.Item(210, 1).Value = 9  # Mercado ID
$ws.Cells.Item(210, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(210, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(210, 4).Value = 44476  # Fecha
$ws.Cells.Item(210, 5).Value = 13  # Codreg
$ws.Cells.Item(210, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(210, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(210, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(210, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(210, 10).Value = 106  # Volumen
$ws.Cells.Item(210, 11).Value = 4000  # Precio mínimo
$ws.Cells.Item(210, 12).Value = 5000  # Precio máximo
$ws.Cells.Item(210, 13).Value = 4500  # Precio promedio ponderado
$ws.Cells.Item(210, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(210, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(210, 16).Value = 450  # Precio $/Kg
$ws.Cells.Item(210, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(210, 18).Value = 'Hortaliza'  # Clasificación

# Row 211
$ws.Cells.Item(211, 1).Value = 9  # Mercado ID
$ws.Cells.Item(211, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(211, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(211, 4).Value = 44468  # Fecha
$ws.Cells.Item(211, 5).Value = 13  # Codreg
$ws.Cells.Item(211, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(211, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(211, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(211, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(211, 10).Value = 250  # Volumen
$ws.Cells.Item(211, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(211, 12).Value = 7000  # Precio máximo
$ws.Cells.Item(211, 13).Value = 6500  # Precio promedio ponderado
$ws.Cells.Item(211, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(211, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(211, 16).Value = 650  # Precio $/Kg
$ws.Cells.Item(211, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(211, 18).Value = 'Hortaliza'  # Clasificación

# Row 212
$ws.Cells.Item(212, 1).Value = 9  # Mercado ID
$ws.Cells.Item(212, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(212, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(212, 4).Value = 44468  # Fecha
$ws.Cells.Item(212, 5).Value = 13  # Codreg
$ws.Cells.Item(212, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(212, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(212, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(212, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(212, 10).Value = 97  # Volumen
$ws.Cells.Item(212, 11).Value = 4000  # Precio mínimo
$ws.Cells.Item(212, 12).Value = 5000  # Precio máximo
$ws.Cells.Item(212, 13).Value = 4495  # Precio promedio ponderado
$ws.Cells.Item(212, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(212, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(212, 16).Value = 450  # Precio $/Kg
$ws.Cells.Item(212, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(212, 18).Value = 'Hortaliza'  # Clasificación

# Row 213
$ws.Cells.Item(213, 1).Value = 9  # Mercado ID
$ws.Cells.Item(213, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(213, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(213, 4).Value = 44386  # Fecha
$ws.Cells.Item(213, 5).Value = 13  # Codreg
$ws.Cells.Item(213, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(213, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(213, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(213, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(213, 10).Value = 340  # Volumen
$ws.Cells.Item(213, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(213, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(213, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(213, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(213, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(213, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(213, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(213, 18).Value = 'Hortaliza'  # Clasificación

# Row 214
$ws.Cells.Item(214, 1).Value = 9  # Mercado ID
$ws.Cells.Item(214, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(214, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(214, 4).Value = 44386  # Fecha
$ws.Cells.Item(214, 5).Value = 13  # Codreg
$ws.Cells.Item(214, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(214, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(214, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(214, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(214, 10).Value = 160  # Volumen
$ws.Cells.Item(214, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(214, 12).Value = 6000  # Precio máximo
$ws.Cells.Item(214, 13).Value = 6000  # Precio promedio ponderado
$ws.Cells.Item(214, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(214, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(214, 16).Value = 600  # Precio $/Kg
$ws.Cells.Item(214, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(214, 18).Value = 'Hortaliza'  # Clasificación

# Row 215
$ws.Cells.Item(215, 1).Value = 9  # Mercado ID
$ws.Cells.Item(215, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(215, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(215, 4).Value = 44384  # Fecha
$ws.Cells.Item(215, 5).Value = 13  # Codreg
$ws.Cells.Item(215, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(215, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(215, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(215, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(215, 10).Value = 52  # Volumen
$ws.Cells.Item(215, 11).Value = 15000  # Precio mínimo
$ws.Cells.Item(215, 12).Value = 16000  # Precio máximo
$ws.Cells.Item(215, 13).Value = 15500  # Precio promedio ponderado
$ws.Cells.Item(215, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(215, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(215, 16).Value = 1550  # Precio $/Kg
$ws.Cells.Item(215, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(215, 18).Value = 'Hortaliza'  # Clasificación

# Row 216
$ws.Cells.Item(216, 1).Value = 9  # Mercado ID
$ws.Cells.Item(216, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(216, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(216, 4).Value = 44384  # Fecha
$ws.Cells.Item(216, 5).Value = 13  # Codreg
$ws.Cells.Item(216, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(216, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(216, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(216, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(216, 10).Value = 25  # Volumen
$ws.Cells.Item(216, 11).Value = 14000  # Precio mínimo
$ws.Cells.Item(216, 12).Value = 14000  # Precio máximo
$ws.Cells.Item(216, 13).Value = 14000  # Precio promedio ponderado
$ws.Cells.Item(216, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(216, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(216, 16).Value = 1400  # Precio $/Kg
$ws.Cells.Item(216, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(216, 18).Value = 'Hortaliza'  # Clasificación

# Row 217
$ws.Cells.Item(217, 1).Value = 9  # Mercado ID
$ws.Cells.Item(217, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(217, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(217, 4).Value = 44363  # Fecha
$ws.Cells.Item(217, 5).Value = 13  # Codreg
$ws.Cells.Item(217, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(217, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(217, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(217, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(217, 10).Value = 250  # Volumen
$ws.Cells.Item(217, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(217, 12).Value = 7000  # Precio máximo
$ws.Cells.Item(217, 13).Value = 6500  # Precio promedio ponderado
$ws.Cells.Item(217, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(217, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(217, 16).Value = 650  # Precio $/Kg
$ws.Cells.Item(217, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(217, 18).Value = 'Hortaliza'  # Clasificación

# Row 218
$ws.Cells.Item(218, 1).Value = 9  # Mercado ID
$ws.Cells.Item(218, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(218, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(218, 4).Value = 44363  # Fecha
$ws.Cells.Item(218, 5).Value = 13  # Codreg
$ws.Cells.Item(218, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(218, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(218, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(218, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(218, 10).Value = 70  # Volumen
$ws.Cells.Item(218, 11).Value = 5000  # Precio mínimo
$ws.Cells.Item(218, 12).Value = 5000  # Precio máximo
$ws.Cells.Item(218, 13).Value = 5000  # Precio promedio ponderado
$ws.Cells.Item(218, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(218, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(218, 16).Value = 500  # Precio $/Kg
$ws.Cells.Item(218, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(218, 18).Value = 'Hortaliza'  # Clasificación

# Row 219
$ws.Cells.Item(219, 1).Value = 9  # Mercado ID
$ws.Cells.Item(219, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(219, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(219, 4).Value = 44306  # Fecha
$ws.Cells.Item(219, 5).Value = 13  # Codreg
$ws.Cells.Item(219, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(219, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(219, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(219, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(219, 10).Value = 340  # Volumen
$ws.Cells.Item(219, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(219, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(219, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(219, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(219, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(219, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(219, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(219, 18).Value = 'Hortaliza'  # Clasificación

# Row 220
$ws.Cells.Item(220, 1).Value = 9  # Mercado ID
$ws.Cells.Item(220, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(220, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(220, 4).Value = 44369  # Fecha
$ws.Cells.Item(220, 5).Value = 13  # Codreg
$ws.Cells.Item(220, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(220, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(220, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(220, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(220, 10).Value = 340  # Volumen
$ws.Cells.Item(220, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(220, 12).Value = 7000  # Precio máximo
$ws.Cells.Item(220, 13).Value = 6500  # Precio promedio ponderado
$ws.Cells.Item(220, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(220, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(220, 16).Value = 650  # Precio $/Kg
$ws.Cells.Item(220, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(220, 18).Value = 'Hortaliza'  # Clasificación

# Row 221
$ws.Cells.Item(221, 1).Value = 9  # Mercado ID
$ws.Cells.Item(221, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(221, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(221, 4).Value = 44301  # Fecha
$ws.Cells.Item(221, 5).Value = 13  # Codreg
$ws.Cells.Item(221, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(221, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(221, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(221, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(221, 10).Value = 340  # Volumen
$ws.Cells.Item(221, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(221, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(221, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(221, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(221, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(221, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(221, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(221, 18).Value = 'Hortaliza'  # Clasificación

# Row 222
$ws.Cells.Item(222, 1).Value = 9  # Mercado ID
$ws.Cells.Item(222, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(222, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(222, 4).Value = 44357  # Fecha
$ws.Cells.Item(222, 5).Value = 13  # Codreg
$ws.Cells.Item(222, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(222, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(222, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(222, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(222, 10).Value = 340  # Volumen
$ws.Cells.Item(222, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(222, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(222, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(222, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(222, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(222, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(222, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(222, 18).Value = 'Hortaliza'  # Clasificación

# Row 223
$ws.Cells.Item(223, 1).Value = 9  # Mercado ID
$ws.Cells.Item(223, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(223, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(223, 4).Value = 44357  # Fecha
$ws.Cells.Item(223, 5).Value = 13  # Codreg
$ws.Cells.Item(223, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(223, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(223, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(223, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(223, 10).Value = 160  # Volumen
$ws.Cells.Item(223, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(223, 12).Value = 6000  # Precio máximo
$ws.Cells.Item(223, 13).Value = 6000  # Precio promedio ponderado
$ws.Cells.Item(223, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(223, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(223, 16).Value = 600  # Precio $/Kg
$ws.Cells.Item(223, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(223, 18).Value = 'Hortaliza'  # Clasificación

# Row 224
$ws.Cells.Item(224, 1).Value = 9  # Mercado ID
$ws.Cells.Item(224, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(224, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(224, 4).Value = 44328  # Fecha
$ws.Cells.Item(224, 5).Value = 13  # Codreg
$ws.Cells.Item(224, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(224, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(224, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(224, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(224, 10).Value = 250  # Volumen
$ws.Cells.Item(224, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(224, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(224, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(224, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(224, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(224, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(224, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(224, 18).Value = 'Hortaliza'  # Clasificación

# Row 225
$ws.Cells.Item(225, 1).Value = 9  # Mercado ID
$ws.Cells.Item(225, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(225, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(225, 4).Value = 44321  # Fecha
$ws.Cells.Item(225, 5).Value = 13  # Codreg
$ws.Cells.Item(225, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(225, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(225, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(225, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(225, 10).Value = 250  # Volumen
$ws.Cells.Item(225, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(225, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(225, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(225, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(225, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(225, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(225, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(225, 18).Value = 'Hortaliza'  # Clasificación

# Row 226
$ws.Cells.Item(226, 1).Value = 9  # Mercado ID
$ws.Cells.Item(226, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(226, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(226, 4).Value = 44223  # Fecha
$ws.Cells.Item(226, 5).Value = 13  # Codreg
$ws.Cells.Item(226, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(226, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(226, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(226, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(226, 10).Value = 120  # Volumen
$ws.Cells.Item(226, 11).Value = 12000  # Precio mínimo
$ws.Cells.Item(226, 12).Value = 14000  # Precio máximo
$ws.Cells.Item(226, 13).Value = 13333  # Precio promedio ponderado
$ws.Cells.Item(226, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(226, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(226, 16).Value = 1333  # Precio $/Kg
$ws.Cells.Item(226, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(226, 18).Value = 'Hortaliza'  # Clasificación

# Row 227
$ws.Cells.Item(227, 1).Value = 9  # Mercado ID
$ws.Cells.Item(227, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(227, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(227, 4).Value = 44298  # Fecha
$ws.Cells.Item(227, 5).Value = 13  # Codreg
$ws.Cells.Item(227, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(227, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(227, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(227, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(227, 10).Value = 250  # Volumen
$ws.Cells.Item(227, 11).Value = 9000  # Precio mínimo
$ws.Cells.Item(227, 12).Value = 9000  # Precio máximo
$ws.Cells.Item(227, 13).Value = 9000  # Precio promedio ponderado
$ws.Cells.Item(227, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(227, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(227, 16).Value = 900  # Precio $/Kg
$ws.Cells.Item(227, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(227, 18).Value = 'Hortaliza'  # Clasificación

# Row 228
$ws.Cells.Item(228, 1).Value = 9  # Mercado ID
$ws.Cells.Item(228, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(228, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(228, 4).Value = 44397  # Fecha
$ws.Cells.Item(228, 5).Value = 13  # Codreg
$ws.Cells.Item(228, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(228, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(228, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(228, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(228, 10).Value = 250  # Volumen
$ws.Cells.Item(228, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(228, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(228, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(228, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(228, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(228, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(228, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(228, 18).Value = 'Hortaliza'  # Clasificación

# Row 229
$ws.Cells.Item(229, 1).Value = 9  # Mercado ID
$ws.Cells.Item(229, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(229, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(229, 4).Value = 44397  # Fecha
$ws.Cells.Item(229, 5).Value = 13  # Codreg
$ws.Cells.Item(229, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(229, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(229, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(229, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(229, 10).Value = 97  # Volumen
$ws.Cells.Item(229, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(229, 12).Value = 6000  # Precio máximo
$ws.Cells.Item(229, 13).Value = 6000  # Precio promedio ponderado
$ws.Cells.Item(229, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(229, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(229, 16).Value = 600  # Precio $/Kg
$ws.Cells.Item(229, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(229, 18).Value = 'Hortaliza'  # Clasificación

# Row 230
$ws.Cells.Item(230, 1).Value = 9  # Mercado ID
$ws.Cells.Item(230, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(230, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(230, 4).Value = 44414  # Fecha
$ws.Cells.Item(230, 5).Value = 13  # Codreg
$ws.Cells.Item(230, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(230, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(230, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(230, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(230, 10).Value = 340  # Volumen
$ws.Cells.Item(230, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(230, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(230, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(230, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(230, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(230, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(230, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(230, 18).Value = 'Hortaliza'  # Clasificación

# Row 231
$ws.Cells.Item(231, 1).Value = 9  # Mercado ID
$ws.Cells.Item(231, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(231, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(231, 4).Value = 44414  # Fecha
$ws.Cells.Item(231, 5).Value = 13  # Codreg
$ws.Cells.Item(231, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(231, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(231, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(231, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(231, 10).Value = 160  # Volumen
$ws.Cells.Item(231, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(231, 12).Value = 6000  # Precio máximo
$ws.Cells.Item(231, 13).Value = 6000  # Precio promedio ponderado
$ws.Cells.Item(231, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(231, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(231, 16).Value = 600  # Precio $/Kg
$ws.Cells.Item(231, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(231, 18).Value = 'Hortaliza'  # Clasificación

# Row 232
$ws.Cells.Item(232, 1).Value = 9  # Mercado ID
$ws.Cells.Item(232, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(232, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(232, 4).Value = 44314  # Fecha
$ws.Cells.Item(232, 5).Value = 13  # Codreg
$ws.Cells.Item(232, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(232, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(232, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(232, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(232, 10).Value = 250  # Volumen
$ws.Cells.Item(232, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(232, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(232, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(232, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(232, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(232, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(232, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(232, 18).Value = 'Hortaliza'  # Clasificación

# Row 233
$ws.Cells.Item(233, 1).Value = 9  # Mercado ID
$ws.Cells.Item(233, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(233, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(233, 4).Value = 44392  # Fecha
$ws.Cells.Item(233, 5).Value = 13  # Codreg
$ws.Cells.Item(233, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(233, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(233, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(233, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(233, 10).Value = 340  # Volumen
$ws.Cells.Item(233, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(233, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(233, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(233, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(233, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(233, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(233, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(233, 18).Value = 'Hortaliza'  # Clasificación

# Row 234
$ws.Cells.Item(234, 1).Value = 9  # Mercado ID
$ws.Cells.Item(234, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(234, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(234, 4).Value = 44392  # Fecha
$ws.Cells.Item(234, 5).Value = 13  # Codreg
$ws.Cells.Item(234, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(234, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(234, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(234, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(234, 10).Value = 160  # Volumen
$ws.Cells.Item(234, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(234, 12).Value = 6000  # Precio máximo
$ws.Cells.Item(234, 13).Value = 6000  # Precio promedio ponderado
$ws.Cells.Item(234, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(234, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(234, 16).Value = 600  # Precio $/Kg
$ws.Cells.Item(234, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(234, 18).Value = 'Hortaliza'  # Clasificación

# Row 235
$ws.Cells.Item(235, 1).Value = 9  # Mercado ID
$ws.Cells.Item(235, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(235, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(235, 4).Value = 44425  # Fecha
$ws.Cells.Item(235, 5).Value = 13  # Codreg
$ws.Cells.Item(235, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(235, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(235, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(235, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(235, 10).Value = 160  # Volumen
$ws.Cells.Item(235, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(235, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(235, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(235, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(235, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(235, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(235, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(235, 18).Value = 'Hortaliza'  # Clasificación

# Row 236
$ws.Cells.Item(236, 1).Value = 9  # Mercado ID
$ws.Cells.Item(236, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(236, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(236, 4).Value = 44425  # Fecha
$ws.Cells.Item(236, 5).Value = 13  # Codreg
$ws.Cells.Item(236, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(236, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(236, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(236, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(236, 10).Value = 97  # Volumen
$ws.Cells.Item(236, 11).Value = 5000  # Precio mínimo
$ws.Cells.Item(236, 12).Value = 5000  # Precio máximo
$ws.Cells.Item(236, 13).Value = 5000  # Precio promedio ponderado
$ws.Cells.Item(236, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(236, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(236, 16).Value = 500  # Precio $/Kg
$ws.Cells.Item(236, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(236, 18).Value = 'Hortaliza'  # Clasificación

# Row 237
$ws.Cells.Item(237, 1).Value = 9  # Mercado ID
$ws.Cells.Item(237, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(237, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(237, 4).Value = 44390  # Fecha
$ws.Cells.Item(237, 5).Value = 13  # Codreg
$ws.Cells.Item(237, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(237, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(237, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(237, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(237, 10).Value = 340  # Volumen
$ws.Cells.Item(237, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(237, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(237, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(237, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(237, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(237, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(237, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(237, 18).Value = 'Hortaliza'  # Clasificación

# Row 238
$ws.Cells.Item(238, 1).Value = 9  # Mercado ID
$ws.Cells.Item(238, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(238, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(238, 4).Value = 44390  # Fecha
$ws.Cells.Item(238, 5).Value = 13  # Codreg
$ws.Cells.Item(238, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(238, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(238, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(238, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(238, 10).Value = 160  # Volumen
$ws.Cells.Item(238, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(238, 12).Value = 6000  # Precio máximo
$ws.Cells.Item(238, 13).Value = 6000  # Precio promedio ponderado
$ws.Cells.Item(238, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(238, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(238, 16).Value = 600  # Precio $/Kg
$ws.Cells.Item(238, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(238, 18).Value = 'Hortaliza'  # Clasificación

# Row 239
$ws.Cells.Item(239, 1).Value = 9  # Mercado ID
$ws.Cells.Item(239, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(239, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(239, 4).Value = 44187  # Fecha
$ws.Cells.Item(239, 5).Value = 13  # Codreg
$ws.Cells.Item(239, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(239, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(239, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(239, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(239, 10).Value = 280  # Volumen
$ws.Cells.Item(239, 11).Value = 8000  # Precio mínimo
$ws.Cells.Item(239, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(239, 13).Value = 8000  # Precio promedio ponderado
$ws.Cells.Item(239, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(239, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(239, 16).Value = 800  # Precio $/Kg
$ws.Cells.Item(239, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(239, 18).Value = 'Hortaliza'  # Clasificación

# Row 240
$ws.Cells.Item(240, 1).Value = 9  # Mercado ID
$ws.Cells.Item(240, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(240, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(240, 4).Value = 44466  # Fecha
$ws.Cells.Item(240, 5).Value = 13  # Codreg
$ws.Cells.Item(240, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(240, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(240, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(240, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(240, 10).Value = 131  # Volumen
$ws.Cells.Item(240, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(240, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(240, 13).Value = 7496  # Precio promedio ponderado
$ws.Cells.Item(240, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(240, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(240, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(240, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(240, 18).Value = 'Hortaliza'  # Clasificación

# Row 241
$ws.Cells.Item(241, 1).Value = 9  # Mercado ID
$ws.Cells.Item(241, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(241, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(241, 4).Value = 44250  # Fecha
$ws.Cells.Item(241, 5).Value = 13  # Codreg
$ws.Cells.Item(241, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(241, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(241, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(241, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(241, 10).Value = 340  # Volumen
$ws.Cells.Item(241, 11).Value = 10000  # Precio mínimo
$ws.Cells.Item(241, 12).Value = 12000  # Precio máximo
$ws.Cells.Item(241, 13).Value = 11000  # Precio promedio ponderado
$ws.Cells.Item(241, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(241, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(241, 16).Value = 1100  # Precio $/Kg
$ws.Cells.Item(241, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(241, 18).Value = 'Hortaliza'  # Clasificación

# Row 242
$ws.Cells.Item(242, 1).Value = 9  # Mercado ID
$ws.Cells.Item(242, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(242, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(242, 4).Value = 44438  # Fecha
$ws.Cells.Item(242, 5).Value = 13  # Codreg
$ws.Cells.Item(242, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(242, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(242, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(242, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(242, 10).Value = 160  # Volumen
$ws.Cells.Item(242, 11).Value = 8000  # Precio mínimo
$ws.Cells.Item(242, 12).Value = 9000  # Precio máximo
$ws.Cells.Item(242, 13).Value = 8500  # Precio promedio ponderado
$ws.Cells.Item(242, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(242, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(242, 16).Value = 850  # Precio $/Kg
$ws.Cells.Item(242, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(242, 18).Value = 'Hortaliza'  # Clasificación

# Row 243
$ws.Cells.Item(243, 1).Value = 9  # Mercado ID
$ws.Cells.Item(243, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(243, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(243, 4).Value = 44201  # Fecha
$ws.Cells.Item(243, 5).Value = 13  # Codreg
$ws.Cells.Item(243, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(243, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(243, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(243, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(243, 10).Value = 280  # Volumen
$ws.Cells.Item(243, 11).Value = 8000  # Precio mínimo
$ws.Cells.Item(243, 12).Value = 10000  # Precio máximo
$ws.Cells.Item(243, 13).Value = 9286  # Precio promedio ponderado
$ws.Cells.Item(243, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(243, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(243, 16).Value = 929  # Precio $/Kg
$ws.Cells.Item(243, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(243, 18).Value = 'Hortaliza'  # Clasificación

# Row 244
$ws.Cells.Item(244, 1).Value = 9  # Mercado ID
$ws.Cells.Item(244, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(244, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(244, 4).Value = 44461  # Fecha
$ws.Cells.Item(244, 5).Value = 13  # Codreg
$ws.Cells.Item(244, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(244, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(244, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(244, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(244, 10).Value = 133  # Volumen
$ws.Cells.Item(244, 11).Value = 6000  # Precio mínimo
$ws.Cells.Item(244, 12).Value = 7000  # Precio máximo
$ws.Cells.Item(244, 13).Value = 6504  # Precio promedio ponderado
$ws.Cells.Item(244, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(244, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(244, 16).Value = 650  # Precio $/Kg
$ws.Cells.Item(244, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(244, 18).Value = 'Hortaliza'  # Clasificación

# Row 245
$ws.Cells.Item(245, 1).Value = 9  # Mercado ID
$ws.Cells.Item(245, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(245, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(245, 4).Value = 44461  # Fecha
$ws.Cells.Item(245, 5).Value = 13  # Codreg
$ws.Cells.Item(245, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(245, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(245, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(245, 9).Value = 'Segunda'  # Calidad
$ws.Cells.Item(245, 10).Value = 61  # Volumen
$ws.Cells.Item(245, 11).Value = 5000  # Precio mínimo
$ws.Cells.Item(245, 12).Value = 5000  # Precio máximo
$ws.Cells.Item(245, 13).Value = 5000  # Precio promedio ponderado
$ws.Cells.Item(245, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(245, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(245, 16).Value = 500  # Precio $/Kg
$ws.Cells.Item(245, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(245, 18).Value = 'Hortaliza'  # Clasificación

# Row 246
$ws.Cells.Item(246, 1).Value = 9  # Mercado ID
$ws.Cells.Item(246, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(246, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(246, 4).Value = 44193  # Fecha
$ws.Cells.Item(246, 5).Value = 13  # Codreg
$ws.Cells.Item(246, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(246, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(246, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(246, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(246, 10).Value = 70  # Volumen
$ws.Cells.Item(246, 11).Value = 8000  # Precio mínimo
$ws.Cells.Item(246, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(246, 13).Value = 8000  # Precio promedio ponderado
$ws.Cells.Item(246, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(246, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(246, 16).Value = 800  # Precio $/Kg
$ws.Cells.Item(246, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(246, 18).Value = 'Hortaliza'  # Clasificación

# Row 247
$ws.Cells.Item(247, 1).Value = 9  # Mercado ID
$ws.Cells.Item(247, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(247, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(247, 4).Value = 44286  # Fecha
$ws.Cells.Item(247, 5).Value = 13  # Codreg
$ws.Cells.Item(247, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(247, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(247, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(247, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(247, 10).Value = 340  # Volumen
$ws.Cells.Item(247, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(247, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(247, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(247, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(247, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(247, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(247, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(247, 18).Value = 'Hortaliza'  # Clasificación

# Row 248
$ws.Cells.Item(248, 1).Value = 9  # Mercado ID
$ws.Cells.Item(248, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(248, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(248, 4).Value = 44389  # Fecha
$ws.Cells.Item(248, 5).Value = 13  # Codreg
$ws.Cells.Item(248, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(248, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(248, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(248, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(248, 10).Value = 250  # Volumen
$ws.Cells.Item(248, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(248, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(248, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(248, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(248, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(248, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(248, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(248, 18).Value = 'Hortaliza'  # Clasificación

# Row 249
$ws.Cells.Item(249, 1).Value = 9  # Mercado ID
$ws.Cells.Item(249, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(249, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(249, 4).Value = 44312  # Fecha
$ws.Cells.Item(249, 5).Value = 13  # Codreg
$ws.Cells.Item(249, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(249, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(249, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(249, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(249, 10).Value = 250  # Volumen
$ws.Cells.Item(249, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(249, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(249, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(249, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(249, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(249, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(249, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(249, 18).Value = 'Hortaliza'  # Clasificación

# Row 250
$ws.Cells.Item(250, 1).Value = 9  # Mercado ID
$ws.Cells.Item(250, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(250, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(250, 4).Value = 44326  # Fecha
$ws.Cells.Item(250, 5).Value = 13  # Codreg
$ws.Cells.Item(250, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(250, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(250, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(250, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(250, 10).Value = 160  # Volumen
$ws.Cells.Item(250, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(250, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(250, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(250, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(250, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(250, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(250, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(250, 18).Value = 'Hortaliza'  # Clasificación
$ws.Cells.Item(250, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Row 251
$ws.Cells.Item(251, 1).Value = 9  # Mercado ID
$ws.Cells.Item(251, 2).Value = 'Vega Central Mapocho de Santiago'  # Mercado
$ws.Cells.Item(251, 3).Value = 'Metropolitana'  # Región
$ws.Cells.Item(251, 4).Value = 44432  # Fecha
$ws.Cells.Item(251, 5).Value = 13  # Codreg
$ws.Cells.Item(251, 6).Value = 100112012  # Categoría ID
$ws.Cells.Item(251, 7).Value = 'Espinaca'  # Categoría
$ws.Cells.Item(251, 8).Value = 'Sin especificar'  # Variedad
$ws.Cells.Item(251, 9).Value = 'Primera'  # Calidad
$ws.Cells.Item(251, 10).Value = 250  # Volumen
$ws.Cells.Item(251, 11).Value = 7000  # Precio mínimo
$ws.Cells.Item(251, 12).Value = 8000  # Precio máximo
$ws.Cells.Item(251, 13).Value = 7500  # Precio promedio ponderado
$ws.Cells.Item(251, 14).Value = '$/cuna 10 kilos'  # Unidad de comercialización
$ws.Cells.Item(251, 15).Value = 'Provincia de Chacabuco'  # Origen
$ws.Cells.Item(251, 16).Value = 750  # Precio $/Kg
$ws.Cells.Item(251, 17).Value = 10  # Kg o Unidades
$ws.Cells.Item(251, 18).Value = 'Hortaliza'  # Clasificación
$ws.Cells.Item(251, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
